$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, copying the header formatting (bold, border, alignment) from A1
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill team record data (Wins/Losses/Ties) for every data row (2-44)
$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 103
    $ws.Cells.Item($r, 30).Value = 59
    $ws.Cells.Item($r, 31).Value = 0
}
